{"js": "// Update the worksheet date and every \"AA\u00d7BB=\" multiplication problem.\n//\n// The document title paragraph holds a date label, and a 5-column table\n// holds the practice problems (data only on every 5th row \u2014 the rows in\n// between are blank spacers). We address every cell positionally\n// (row/column index) rather than by searching for the old text, since one\n// of the \"old\" values re-appears as a \"new\" value elsewhere in the table\n// (table row 0 \"42\u00d748=\" becomes \"44\u00d752=\", while table row 19 already holds\n// \"44\u00d752=\" and becomes \"79\u00d756=\"), and blind text search-and-replace run in\n// the wrong order could clobber the wrong cell.\n\nconst body = context.document.body;\n\n// --- 1. Title paragraph: \"2024-12-28 Saturday\" -> \"2024-12-29 Sunday\" ---\nconst title = body.paragraphs.getFirst();\ntitle.insertText(\"2024-12-29 Sunday\", \"Replace\");\n\n// --- 2. Multiplication table ---\nconst table = body.tables.getFirst();\n\n// Table row index (0-based) -> new values for that row's 5 columns.\nconst newRowValues = {\n  0: [\"44\u00d752=\", \"61\u00d727=\", \"12\u00d797=\", \"47\u00d758=\", \"89\u00d762=\"],\n  4: [\"38\u00d789=\", \"43\u00d774=\", \"93\u00d765=\", \"40\u00d761=\", \"36\u00d714=\"],\n  9: [\"16\u00d714=\", \"20\u00d778=\", \"78\u00d734=\", \"37\u00d717=\", \"56\u00d773=\"],\n  14: [\"41\u00d748=\", \"57\u00d784=\", \"35\u00d786=\", \"53\u00d780=\", \"14\u00d735=\"],\n  19: [\"79\u00d756=\", \"17\u00d762=\", \"31\u00d720=\", \"64\u00d745=\", \"14\u00d794=\"]\n};\n\nfor (const rowIndex of Object.keys(newRowValues)) {\n  const r = Number(rowIndex);\n  const values = newRowValues[r];\n  for (let c = 0; c < values.length; c++) {\n    table.getCell(r, c).value = values[c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date and every \"AA\u00d7BB=\" multiplication problem.\n#\n# The document title paragraph holds a date label, and a 5-column table\n# holds the practice problems (data only on every 5th row \u2014 the rows in\n# between are blank spacers). Cells are addressed positionally\n# (1-based row/column index) rather than via Find/Replace on the old text,\n# since a couple of the \"old\" values re-appear as \"new\" values elsewhere in\n# the table (e.g. cell (1,1) \"42\u00d748=\" becomes \"44\u00d752=\", while the cell that\n# currently holds \"44\u00d752=\" becomes \"79\u00d756=\"), and a naive text\n# search-and-replace run in the wrong order could clobber the wrong cell.\n\n$d = $word.ActiveDocument\n\n# --- 1. Title paragraph: \"2024-12-28 Saturday\" -> \"2024-12-29 Sunday\" ---\n$d.Paragraphs.Item(1).Range.Text = \"2024-12-29 Sunday\"\n\n# --- 2. Multiplication table ---\n$t = $d.Tables.Item(1)\n\n# 1-based table row -> new values for that row's 5 columns.\n$newRowValues = @{\n    1  = @(\"44\u00d752=\", \"61\u00d727=\", \"12\u00d797=\", \"47\u00d758=\", \"89\u00d762=\")\n    5  = @(\"38\u00d789=\", \"43\u00d774=\", \"93\u00d765=\", \"40\u00d761=\", \"36\u00d714=\")\n    10 = @(\"16\u00d714=\", \"20\u00d778=\", \"78\u00d734=\", \"37\u00d717=\", \"56\u00d773=\")\n    15 = @(\"41\u00d748=\", \"57\u00d784=\", \"35\u00d786=\", \"53\u00d780=\", \"14\u00d735=\")\n    20 = @(\"79\u00d756=\", \"17\u00d762=\", \"31\u00d720=\", \"64\u00d745=\", \"14\u00d794=\")\n}\n\nforeach ($row in $newRowValues.Keys) {\n    $values = $newRowValues[$row]\n    for ($col = 1; $col -le $values.Length; $col++) {\n        $t.Cell($row, $col).Range.Text = $values[$col - 1]\n    }\n}\n"}
